$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest crypto data.
# Price cells keep their original text format (dotted thousands/decimal
# separators), so force the cell to Text before writing the value; this
# prevents Excel from auto-converting strings like "1.003" into a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.449.58"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.568.39"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.13"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3726"
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.29"
$ws.Range("E8").Value = "  -3.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3314"
$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07475"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.73"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.962"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.910"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.06"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001114"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.03"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06766"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.365"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.06"
$ws.Range("E23").Value = "  -1.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.438.79"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.394"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.558"
$ws.Range("E26").Value = "  -4.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.35"
$ws.Range("E27").Value = "  +2.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.66"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.018"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.78"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.747.54"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.051"
$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.126"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.619"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08289"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02450"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2267"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06393"
$ws.Range("E39").Value = "  -2.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.364"
$ws.Range("E40").Value = "  -1.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.291"
$ws.Range("E41").Value = "  -4.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6275"
$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.83"
$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6135"
$ws.Range("E46").Value = "  +4.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.784"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.042"
$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.22"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.212"
$ws.Range("E50").Value = "  -2.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07245"
$ws.Range("E51").Value = "  -0.83%  "
